$wb = $excel.ActiveWorkbook

# --- Update "Priority" column (E) to "ht" for specific rows on zh-cn and de-de sheets ---
$rows = @(8, 9, 10, 11, 13, 14)

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($r in $rows) {
        $ws.Cells.Item($r, 5).Value = "ht"
    }
}

# --- Update "Latest Handoff Datetime" timestamps ---

# zh-cn sheet (column H): 2016-08-27 20:21:44 -> 2016-08-27 20:21:59
$wsZh = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZh.Cells.Item($r, 8).Value = "2016-08-27 20:21:59"
}

# de-de sheet (column H): 2016-08-27 20:21:49 -> 2016-08-27 20:22:09
$wsDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDe.Cells.Item($r, 8).Value = "2016-08-27 20:22:09"
}

# Overview sheet (column G, "Latest HO Xliff Generate Date"): same update as de-de
$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-27 20:22:09"
}
